$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-4: only column H (CV Score) values change ---
$ws.Range("H2").Value = 4.352897191762352
$ws.Range("H3").Value = 4.334840652231079
$ws.Range("H4").Value = 4.333967895825444

# --- Rows 5 and 6: CatBoostRegressor and XGBRegressor swap places,
#     each with updated hyperparameters/metrics ---
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "CatBoostRegressor"
$ws.Range("C5").Value = "{'depth': 3, 'iterations': 1000, 'l2_leaf_reg': 1, 'learning_rate': 0.01}"
$ws.Range("D5").Value = 0.8742270774403477
$ws.Range("E5").Value = 4.279435844213173
$ws.Range("F5").Value = 5.532213095949007
$ws.Range("G5").Value = 30.6053817389897
$ws.Range("H5").Value = 4.456699827494

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "XGBRegressor"
$ws.Range("C6").Value = "{'learning_rate': 0.1, 'max_depth': 3, 'n_estimators': 100, 'random_state': 42, 'subsample': 0.8}"
$ws.Range("D6").Value = 0.869636058807373
$ws.Range("E6").Value = 4.366933994293213
$ws.Range("F6").Value = 5.632278085654206
$ws.Range("G6").Value = 31.72255643414061
$ws.Range("H6").Value = 4.624904854297638

# --- Row 7: RandomForestRegressor hyperparameters and metrics updated ---
$ws.Range("C7").Value = "{'max_depth': None, 'min_samples_leaf': 5, 'min_samples_split': 2, 'n_estimators': 50, 'random_state': 42}"
$ws.Range("D7").Value = 0.848044588094985
$ws.Range("E7").Value = 4.652564932331449
$ws.Range("F7").Value = 6.080837672912122
$ws.Range("G7").Value = 36.97658680430732
$ws.Range("H7").Value = 4.777222735225072

# --- Rows 8-10: only column H (CV Score) values change ---
$ws.Range("H8").Value = 4.98794190290956
$ws.Range("H9").Value = 5.293914599011564
$ws.Range("H10").Value = 6.853159430414824
